$d = $word.ActiveDocument

# Find the paragraph that contains "SPECIFIED_HEAD:" with placeholder "(XXX--IDD_XXX)"
# and rewrite its trailing "(XXX--IDD_XXX)" portion with the new content.

foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    if ($r.Text -like "SPECIFIED_HEAD:*" -and $r.Text -like "*XXX*") {
        # Find the run of text "(XXX--IDD_XXX)" inside this paragraph and replace piecewise
        $findRange = $r.Duplicate
        $findRange.Find.Execute("(XXX--IDD_XXX)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
        if ($findRange.Find.Found) {
            # findRange now covers "(XXX--IDD_XXX)"
            $findRange.Text = "(CSpecifiedHeadPropsPage -- IDD_PROPS_SPECIFIED_HEAD2)"
        }
        break
    }
}

Write-Host "done"
